# Add new device/password mapping rows (mongoDB support) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newIPs = @(
    "10.127.125.221",
    "10.127.125.221",
    "10.127.125.220",
    "10.127.125.238",
    "10.127.125.222"
)

$startRow = 25
for ($i = 0; $i -lt $newIPs.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newIPs[$i]
    $ws.Cells.Item($row, 2).Value = "nbv_12345"
    $ws.Cells.Item($row, 1).Font.Color = 0
    $ws.Cells.Item($row, 2).Font.Color = 0
}

$ws.Range("F24").Select()
